# Update the simulated-game transition-probability matrix on Sheet1 with
# refreshed values (more games simulated -> slightly different probabilities).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5811320754716981
$ws.Range("J2").Value = 0.003773584905660377
$ws.Range("P2").Value = 0.1471698113207547
$ws.Range("S2").Value = 0.06792452830188679
$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.03125
$ws.Range("P3").Value = 0.74375
$ws.Range("S3").Value = 0.21875
$ws.Range("P4").Value = 0.5686274509803921
$ws.Range("S4").Value = 0.4313725490196079
$ws.Range("B6").Value = 0.04205607476635514
$ws.Range("D6").Value = 0.01401869158878505
$ws.Range("F6").Value = 0.04205607476635514
$ws.Range("J6").Value = 0.2616822429906542
$ws.Range("O6").Value = 0.01401869158878505
$ws.Range("Q6").Value = 0.1635514018691589
$ws.Range("R6").Value = 0.04672897196261682
$ws.Range("S6").Value = 0.4158878504672897
$ws.Range("B7").Value = 0.09395973154362416
$ws.Range("D7").Value = 0.01342281879194631
$ws.Range("F7").Value = 0.04026845637583892
$ws.Range("J7").Value = 0.0738255033557047
$ws.Range("O7").Value = 0.02684563758389262
$ws.Range("Q7").Value = 0.2348993288590604
$ws.Range("R7").Value = 0.06711409395973154
$ws.Range("S7").Value = 0.4496644295302014
$ws.Range("B8").Value = 0.09864603481624758
$ws.Range("D8").Value = 0.03094777562862669
$ws.Range("E8").Value = 0.003868471953578337
$ws.Range("F8").Value = 0.06189555125725339
$ws.Range("J8").Value = 0.0735009671179884
$ws.Range("O8").Value = 0.02127659574468085
$ws.Range("Q8").Value = 0.1992263056092843
$ws.Range("R8").Value = 0.07156673114119923
$ws.Range("S8").Value = 0.4390715667311412
$ws.Range("B9").Value = 0.08571428571428572
$ws.Range("D9").Value = 0.01904761904761905
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.07142857142857142
$ws.Range("O9").Value = 0.03333333333333333
$ws.Range("Q9").Value = 0.1952380952380952
$ws.Range("R9").Value = 0.08571428571428572
$ws.Range("S9").Value = 0.4428571428571428
$ws.Range("B10").Value = 0.1052631578947368
$ws.Range("D10").Value = 0.02497769848349688
$ws.Range("E10").Value = 0.0008920606601248885
$ws.Range("F10").Value = 0.07136485280999108
$ws.Range("J10").Value = 0.07760927743086529
$ws.Range("O10").Value = 0.01427297056199822
$ws.Range("Q10").Value = 0.215878679750223
$ws.Range("R10").Value = 0.07760927743086529
$ws.Range("S10").Value = 0.4121320249776985
$ws.Range("G11").Value = 0.1374407582938389
$ws.Range("J11").Value = 0.06635071090047394
$ws.Range("K11").Value = 0.1658767772511848
$ws.Range("L11").Value = 0.6161137440758294
$ws.Range("S11").Value = 0.01421800947867299
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.1333333333333333
$ws.Range("K12").Value = 0.01481481481481482
$ws.Range("L12").Value = 0.02962962962962963
$ws.Range("S12").Value = 0.02222222222222222
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.3461538461538461
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("F15").Value = 0.01421800947867299
$ws.Range("H15").Value = 0.1611374407582938
$ws.Range("I15").Value = 0.0995260663507109
$ws.Range("J15").Value = 0.3080568720379147
$ws.Range("K15").Value = 0.07109004739336493
$ws.Range("M15").Value = 0.01421800947867299
$ws.Range("N15").Value = 0.004739336492890996
$ws.Range("O15").Value = 0.02843601895734597
$ws.Range("S15").Value = 0.2985781990521327
$ws.Range("F16").Value = 0.03333333333333333
$ws.Range("H16").Value = 0.2055555555555555
$ws.Range("I16").Value = 0.07222222222222222
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.01666666666666667
$ws.Range("O16").Value = 0.06666666666666667
$ws.Range("S16").Value = 0.1222222222222222
$ws.Range("F17").Value = 0.03311258278145696
$ws.Range("H17").Value = 0.1832229580573951
$ws.Range("I17").Value = 0.1059602649006623
$ws.Range("J17").Value = 0.4216335540838852
$ws.Range("K17").Value = 0.0772626931567329
$ws.Range("M17").Value = 0.01324503311258278
$ws.Range("O17").Value = 0.06181015452538632
$ws.Range("S17").Value = 0.1037527593818985
$ws.Range("F18").Value = 0.01234567901234568
$ws.Range("H18").Value = 0.1851851851851852
$ws.Range("I18").Value = 0.1172839506172839
$ws.Range("J18").Value = 0.4506172839506173
$ws.Range("K18").Value = 0.04320987654320987
$ws.Range("O18").Value = 0.06790123456790123
$ws.Range("S18").Value = 0.1234567901234568
$ws.Range("F19").Value = 0.01387818041634541
$ws.Range("H19").Value = 0.2613723978411719
$ws.Range("I19").Value = 0.08558211256746338
$ws.Range("J19").Value = 0.3770239013107171
$ws.Range("K19").Value = 0.07710100231303008
$ws.Range("M19").Value = 0.01079414032382421
$ws.Range("N19").Value = 0.0007710100231303007
$ws.Range("O19").Value = 0.06784888203546646
$ws.Range("S19").Value = 0.1056283731688512
